# Slide 1: "Tap'n Cook" logo textbox.
#   - "Tap'n"  -> "Tap'N"   (capitalize the N)
#   - " Cook"  -> split into two runs: " " and "Cook" (same formatting)
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 6")
$tr = $sh.TextFrame.TextRange

# "Tap'n Cook" -> characters 1-5 are "Tap'n" (using the curly
# right-single-quote, U+2019, as used in the original text).
$apos = [char]0x2019
$run1 = $tr.Characters(1, 5)
$run1.Text = "Tap" + $apos + "N"

# Characters 7-10 are "Cook" (character 6 is the leading space). Re-assigning
# this sub-range splits the existing " Cook" run into a " " run followed by
# a new "Cook" run, matching the target markup.
$run2 = $tr.Characters(7, 4)
$run2.Text = "Cook"
